$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet previously ran through row 343 (9 Aug 2021). This update appends
# daily figures through 23 Aug 2021 (rows 344-357), per "aggiornamento al
# 23 agosto 2021".
#
# First, stamp the new rows (A344:D357) with the same formatting as the last
# existing data row (343) - this carries over column A's date-stamped style
# (centered, bordered, YYYY-MM-DD HH:MM:SS number format) to the new cells.
$ws.Range("A343:D343").Copy($ws.Range("A344:D357"))

# Now fill in the actual values for each new day.
$newData = @(
    @(344, 44418, 0, 3, 91.6030534351145),
    @(345, 44419, 0, 3, 91.6030534351145),
    @(346, 44420, 1, 3, 91.6030534351145),
    @(347, 44421, 1, 4, 122.1374045801527),
    @(348, 44422, 0, 3, 91.6030534351145),
    @(349, 44423, 0, 3, 91.6030534351145),
    @(350, 44424, 0, 2, 61.06870229007634),
    @(351, 44425, 0, 2, 61.06870229007634),
    @(352, 44426, 0, 2, 61.06870229007634),
    @(353, 44427, 0, 1, 30.53435114503817),
    @(354, 44428, 0, 0, 0),
    @(355, 44429, 0, 0, 0),
    @(356, 44430, 0, 0, 0),
    @(357, 44431, 0, 0, 0)
)

foreach ($entry in $newData) {
    $r = $entry[0]
    $ws.Cells.Item($r, 1).Value = $entry[1]
    $ws.Cells.Item($r, 2).Value = $entry[2]
    $ws.Cells.Item($r, 3).Value = $entry[3]
    $ws.Cells.Item($r, 4).Value = $entry[4]
}
